# Update generated-output numbers for "丽水·AEO纯白礼赞动漫嘉年华" row (row 2)
# on both the "展览" and "全部类型" sheets:
#   F2 (想去人数 / want-to-go count): 1053 -> 1056
#   G2 (最低票价 / lowest price): 65 -> "已售罄" (Sold out)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1056
    $ws.Range("G2").Value = "已售罄"
}
